$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.668.51"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.246.16"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -3.10%  "
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.76%  "
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -8.13%  "
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.37%  "
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.575"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.94%  "
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.242.83"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -3.01%  "
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -6.99%  "
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.569"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.87%  "
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "44.68"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -5.54%  "
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.95%  "
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "653.81"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.26%  "
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.781.21"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.59%  "
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.22"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.83%  "
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.765.34"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.25%  "
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.30%  "
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.254.57"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.81%  "
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -4.47%  "
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -4.37%  "
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.875"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -3.35%  "
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -5.09%  "
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.27"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.16%  "
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "96.66"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -4.17%  "
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -4.48%  "
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -6.46%  "
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.04"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -6.69%  "
$ws.Range("E28").Style = "Normal"

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "32.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.35%  "
$ws.Range("E29").Style = "Normal"

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.25"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -4.62%  "
$ws.Range("E30").Style = "Normal"

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.80"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.09%  "
$ws.Range("E31").Style = "Normal"

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "568.66"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -6.45%  "
$ws.Range("E32").Style = "Normal"

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "10.85"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.10%  "
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.776.71"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.81%  "
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.18%  "
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -4.01%  "
$ws.Range("E36").Style = "Normal"

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.38"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -12.45%  "
$ws.Range("E37").Style = "Normal"

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "55.15"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.18%  "
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.27%  "
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "InjectiveProtocol"
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "32.04"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.75%  "
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "Fetch.AI"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.58"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -6.74%  "
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "Stacks"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.99"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -7.17%  "
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "PEPE"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0₃0650"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -7.50%  "
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.21"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -4.49%  "
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.323"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -5.98%  "
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0396"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -5.56%  "
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.11%  "
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.04%  "
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.53"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.57%  "
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.59%  "
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "128.38"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.77%  "
$ws.Range("E51").Style = "Normal"

